$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "66.259.75"
$ws.Range("E2").Value = "  +5.20%  "
# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.799.96"
$ws.Range("E3").Value = "  +8.11%  "
# Row 4
$ws.Range("E4").Value = "  -0.27%  "
# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "428.62"
$ws.Range("E5").Value = "  +10.34%  "
# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "139.79"
$ws.Range("E6").Value = "  +15.10%  "
# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.615"
$ws.Range("E7").Value = "  +6.31%  "
# Row 8
$ws.Range("E8").Value = "  -0.02%  "
# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.744"
$ws.Range("E9").Value = "  +11.20%  "
# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.157"
$ws.Range("E10").Value = "  +5.33%  "
# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0000327"
$ws.Range("E11").Value = "  +4.14%  "
# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "43.43"
$ws.Range("E12").Value = "  +13.63%  "
# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "10.68"
# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.387.91"
$ws.Range("E14").Value = "  +7.64%  "
# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.10"
$ws.Range("E15").Value = "  +17.02%  "
# Row 16
$ws.Range("E16").Value = "  +1.49%  "
# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.779.64"
$ws.Range("E17").Value = "  +7.17%  "
# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "20.20"
$ws.Range("E18").Value = "  +9.35%  "
# Row 19
$ws.Range("E19").Value = "  +13.15%  "
# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "66.355.48"
$ws.Range("E20").Value = "  +5.30%  "
# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "413.05"
$ws.Range("E21").Value = "  +6.45%  "
# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "15.26"
$ws.Range("E22").Value = "  +11.17%  "
# Row 23
$ws.Range("E23").Value = "  +15.80%  "
# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "85.86"
$ws.Range("E24").Value = "  +6.39%  "
# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "37.21"
$ws.Range("E25").Value = "  +11.96%  "
# Row 26
$ws.Range("B26").Value = "RenderToken"
$ws.Range("C26").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.79"
$ws.Range("E26").Value = "  +47.57%  "
# Row 27
$ws.Range("B27").Value = "PancakeSwap"
$ws.Range("C27").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "3.30"
$ws.Range("E27").Value = "  +12.30%  "
# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.88"
$ws.Range("E28").Value = "  +15.67%  "
# Row 29
$ws.Range("E29").Value = "  -0.80%  "
# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "13.92"
$ws.Range("E30").Value = "  +20.15%  "
# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "707.34"
$ws.Range("E31").Value = "  +7.51%  "
# Row 32
$ws.Range("E32").Value = "  +19.39%  "
# Row 33
$ws.Range("E33").Value = "  +8.65%  "
# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "40.15"
$ws.Range("E34").Value = "  +11.44%  "
# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.88"
$ws.Range("E35").Value = "  +45.73%  "
# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.00"
$ws.Range("E36").Value = "  +0.04%  "
# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "56.45"
$ws.Range("E38").Value = "  +6.70%  "
# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0477"
$ws.Range("E39").Value = "  +11.00%  "
# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.64"
$ws.Range("E40").Value = "  +54.87%  "
# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0₃0689"
$ws.Range("E41").Value = "  +11.18%  "
# Row 42
$ws.Range("E42").Value = "  +9.59%  "
# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.85"
$ws.Range("E43").Value = "  +8.30%  "
# Row 44
$ws.Range("E44").Value = "  +0.11%  "
# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.36"
$ws.Range("E45").Value = "  +11.38%  "
# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.322"
$ws.Range("E46").Value = "  +18.98%  "
# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.16"
$ws.Range("E47").Value = "  +3.99%  "
# Row 48
$ws.Range("E48").Value = "  +7.95%  "
# Row 49
$ws.Range("E49").Value = "  +8.30%  "
# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "142.78"
$ws.Range("E50").Value = "  +3.39%  "
# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.82"
$ws.Range("E51").Value = "  +8.51%  "
